$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new row at the top of the data block (row 81),
# pushing the existing historical rows (81-180) down to (82-181), and
# populate the new row with this week's price report for Ají.
$ws.Rows.Item(81).Insert()

$ws.Range("A81").Value = 4
$ws.Range("B81").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C81").Value = "Los Lagos"
$ws.Range("D81").Value = 44539
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 100112021
$ws.Range("G81").Value = "Ají"
$ws.Range("H81").Value = "Inferno"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 70
$ws.Range("K81").Value = 18000
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = 19143
$ws.Range("N81").Value = "$/caja 12 kilos"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 1595
$ws.Range("Q81").Value = 12
$ws.Range("R81").Value = "Hortaliza"
